$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.4
$ws.Range("J2").Value = 1.07
$ws.Range("K2").Value = 9
$ws.Range("P2").Value = 1.44
$ws.Range("Q2").Value = 2.63
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 1.53
$ws.Range("U2").Value = 6
$ws.Range("AG2").Value = 29
$ws.Range("AI2").Value = 81

# Row 3
$ws.Range("K3").Value = 8.5
$ws.Range("W3").Value = 26
$ws.Range("Y3").Value = 34

# Row 4
$ws.Range("G4").Value = 2.3
$ws.Range("I4").Value = 3.1
$ws.Range("N4").Value = 2.03
$ws.Range("O4").Value = 1.87
$ws.Range("T4").Value = 8.5
$ws.Range("U4").Value = 12
$ws.Range("W4").Value = 23
$ws.Range("Y4").Value = 34
$ws.Range("Z4").Value = 10
$ws.Range("AD4").Value = 301
$ws.Range("AE4").Value = 11
$ws.Range("AI4").Value = 29

# Row 6
$ws.Range("G6").Value = 2
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 3.7
$ws.Range("R6").Value = 1.62
$ws.Range("V6").Value = 8.5
$ws.Range("Y6").Value = 21
$ws.Range("Z6").Value = 13
$ws.Range("AA6").Value = 7
$ws.Range("AF6").Value = 21
$ws.Range("AG6").Value = 13
$ws.Range("AH6").Value = 41

# Row 7
$ws.Range("H7").Value = 4
$ws.Range("I7").Value = 4.75
$ws.Range("L7").Value = 1.18
$ws.Range("M7").Value = 4.5
$ws.Range("R7").Value = 1.62
$ws.Range("S7").Value = 2.2
$ws.Range("AG7").Value = 15

# Row 8
$ws.Range("L8").Value = 1.57
$ws.Range("M8").Value = 2.25
$ws.Range("R8").Value = 2.25
$ws.Range("S8").Value = 1.57

# Row 9
$ws.Range("S9").Value = 1.62

# Row 10
$ws.Range("K10").Value = 12
$ws.Range("N10").Value = 1.75
$ws.Range("O10").Value = 2.05
$ws.Range("P10").Value = 1.33
$ws.Range("Q10").Value = 3.25
$ws.Range("AB10").Value = 13
$ws.Range("AJ10").Value = 34

# Row 13
$ws.Range("G13").Value = 4.75
$ws.Range("H13").Value = 3.4
$ws.Range("I13").Value = 1.75
$ws.Range("T13").Value = 11
$ws.Range("U13").Value = 23
$ws.Range("AH13").Value = 13

# Row 14
$ws.Range("G14").Value = 2.9
$ws.Range("I14").Value = 2.2
$ws.Range("L14").Value = 1.17
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = 1.6
$ws.Range("O14").Value = 2.3
$ws.Range("R14").Value = 1.5
$ws.Range("S14").Value = 2.5
$ws.Range("W14").Value = 34
$ws.Range("AB14").Value = 11
$ws.Range("AG14").Value = 9
$ws.Range("AH14").Value = 21
$ws.Range("AI14").Value = 15

# Row 15
$ws.Range("G15").Value = 2.35
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 8.5
$ws.Range("P15").Value = 1.44
$ws.Range("Q15").Value = 2.63
$ws.Range("U15").Value = 11
$ws.Range("V15").Value = 9.5
$ws.Range("W15").Value = 23
$ws.Range("AE15").Value = 8.5
$ws.Range("AF15").Value = 15
$ws.Range("AH15").Value = 34
$ws.Range("AI15").Value = 26

# Row 17
$ws.Range("K17").Value = 10

# Row 21
$ws.Range("G21").Value = 7.5
$ws.Range("H21").Value = 4.75
$ws.Range("I21").Value = 1.36
$ws.Range("K21").Value = 12
$ws.Range("N21").Value = 1.73
$ws.Range("O21").Value = 2.08
$ws.Range("P21").Value = 1.33
$ws.Range("Q21").Value = 3.25
$ws.Range("U21").Value = 41
$ws.Range("AA21").Value = 9.5
$ws.Range("AE21").Value = 7

# Row 22
$ws.Range("N22").Value = 1.73
$ws.Range("O22").Value = 2.08

# Row 26
$ws.Range("L26").Value = 1.29
$ws.Range("M26").Value = 3.5
$ws.Range("N26").Value = 1.98
$ws.Range("O26").Value = 1.88

# Row 28
$ws.Range("J28").Value = 1.05
$ws.Range("K28").Value = 11
$ws.Range("N28").Value = 1.88
$ws.Range("O28").Value = 1.98
